$d = $word.ActiveDocument

# Locate the "LOQ4057 - ... (Requisito fraco)" line (including its trailing line break)
# inside the Requisitos list and cut it out.
$src = $d.Content
$found = $src.Find.Execute(
    "LOQ4057 -  Operações Unitárias III  (Requisito fraco)^l",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if (-not $found) {
    throw "Could not find the LOQ4057 requisito line"
}

$movedText = $src.Text
$src.Delete()

# Locate the start of the "LOB1056 - ..." line and insert the moved text right before it,
# so LOQ4057 now appears first in the Requisitos list.
$dest = $d.Content
$found2 = $dest.Find.Execute(
    "LOB1056 -  Introdução aos Métodos Numéricos e Computacionais  (Requisito fraco)",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if (-not $found2) {
    throw "Could not find the LOB1056 requisito line"
}

$target = $d.Range($dest.Start, $dest.Start)
$target.InsertBefore($movedText)
